# Commit: "fix for tc and implement bidirectional relationship among all entities"
#
# On the "Future Road Map" sheet, add a new row below the existing list of
# items ("Calling 1 tc in another tc" is the last one in B7) with the new
# test case description, then leave the selection on the next empty cell
# below it (as Excel does after typing a value into a cell and pressing
# Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Future Road Map")
$ws.Activate()

$ws.Range("B8").Value = "In get modules test cases array displaying twice"

$ws.Range("B9").Select()
